# BSD - Protokoll: add two new entries to the "Sprint 3" tracking table.
#
# The sprint table is laid out three times across the sheet (columns C-G,
# M-Q, W-AA), each a parallel "Aufgabe / Datum / Von / Bis / Status" block.
# Row 35 gains a "User Management" (left block) and a "Google Maps routing
# recherchiert" (middle block) entry; row 30's right-hand block (previously
# blank) gains a "Backend Security" entry.
#
# New shared strings are appended in the order they are first referenced,
# so write the row 35 cells (User Management, then Google Maps routing
# recherchiert) before the row 30 cell (Backend Security) to land them at
# shared-string indices 56, 57, 58 respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35 - left block: User Management
$ws.Range("C35").Value = "User Management"
$ws.Range("D35").Value = 42802
$ws.Range("E35").Value = 0.59027777777777779
$ws.Range("F35").Value = 0.65972222222222221
$ws.Range("G35").Value = 0.3

# Row 35 - middle block: Google Maps routing recherchiert
$ws.Range("M35").Value = "Google Maps routing recherchiert"
$ws.Range("N35").Value = 42802
$ws.Range("O35").Value = 0.59027777777777779
$ws.Range("P35").Value = 0.65972222222222221
$ws.Range("Q35").Value = 1

# Row 30 - right block: Backend Security
$ws.Range("W30").Value = "Backend Security"
$ws.Range("X30").Value = 42802
$ws.Range("Y30").Value = 0.59027777777777779
$ws.Range("Z30").Value = 0.65972222222222221
$ws.Range("AA30").Value = 1

# Update the saved cursor/selection position to match where editing ended.
$ws.Range("K33").Select()
